$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.171993374824524
$ws.Range("B1").Value = 2.330319881439209
$ws.Range("C1").Value = 3.287084817886353
$ws.Range("D1").Value = 1.482364773750305
$ws.Range("E1").Value = 1.184820413589478
